# Natmi following Dr Hou advice
# Expand the Ptprc-Cd22 L-R pair results to include the "sCs" (stromal/other) cluster
# alongside the existing ECs / FAPs / M2 clusters, for both sending and target sides.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ptprc"
$ws.Cells.Item(2, 3).Value = "Cd22"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 205.313027
$ws.Cells.Item(2, 8).Value = 615.9390810000001
$ws.Cells.Item(2, 9).Value = 0.435242422384838
$ws.Cells.Item(2, 10).Value = 0.435242422384838
$ws.Cells.Item(2, 11).Value = 1.0
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.5751623333333333
$ws.Cells.Item(2, 14).Value = 1.725487
$ws.Cells.Item(2, 15).Value = 0.08167574423707133
$ws.Cells.Item(2, 16).Value = 0.08167574423707133
$ws.Cells.Item(2, 17).Value = 118.0883196730497
$ws.Cells.Item(2, 18).Value = 1062.794877057447
$ws.Cells.Item(2, 19).Value = 0.0355487487718274
$ws.Cells.Item(2, 20).Value = 0.0355487487718274

# Row 3: ECs -> M2
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ptprc"
$ws.Cells.Item(3, 3).Value = "Cd22"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 205.313027
$ws.Cells.Item(3, 8).Value = 615.9390810000001
$ws.Cells.Item(3, 9).Value = 0.435242422384838
$ws.Cells.Item(3, 10).Value = 0.435242422384838
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 6.466858999999999
$ws.Cells.Item(3, 14).Value = 19.400577
$ws.Cells.Item(3, 15).Value = 0.9183242557629286
$ws.Cells.Item(3, 16).Value = 0.9183242557629286
$ws.Cells.Item(3, 17).Value = 1327.730396472193
$ws.Cells.Item(3, 18).Value = 11949.57356824974
$ws.Cells.Item(3, 19).Value = 0.3996936736130106
$ws.Cells.Item(3, 20).Value = 0.3996936736130106

# Row 4: FAPs -> ECs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ptprc"
$ws.Cells.Item(4, 3).Value = "Cd22"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 0.186821
$ws.Cells.Item(4, 8).Value = 0.5604629999999999
$ws.Cells.Item(4, 9).Value = 0.0003960412341120362
$ws.Cells.Item(4, 10).Value = 0.0003960412341120362
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.5751623333333333
$ws.Cells.Item(4, 14).Value = 1.725487
$ws.Cells.Item(4, 15).Value = 0.08167574423707133
$ws.Cells.Item(4, 16).Value = 0.08167574423707133
$ws.Cells.Item(4, 17).Value = 0.1074524022756667
$ws.Cells.Item(4, 18).Value = 0.9670716204809999
$ws.Cells.Item(4, 19).Value = 0.00003234696254466875
$ws.Cells.Item(4, 20).Value = 0.00003234696254466875

# Row 5: FAPs -> M2
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ptprc"
$ws.Cells.Item(5, 3).Value = "Cd22"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 0.186821
$ws.Cells.Item(5, 8).Value = 0.5604629999999999
$ws.Cells.Item(5, 9).Value = 0.0003960412341120362
$ws.Cells.Item(5, 10).Value = 0.0003960412341120362
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 6.466858999999999
$ws.Cells.Item(5, 14).Value = 19.400577
$ws.Cells.Item(5, 15).Value = 0.9183242557629286
$ws.Cells.Item(5, 16).Value = 0.9183242557629286
$ws.Cells.Item(5, 17).Value = 1.208145065239
$ws.Cells.Item(5, 18).Value = 10.873305587151
$ws.Cells.Item(5, 19).Value = 0.0003636942715673674
$ws.Cells.Item(5, 20).Value = 0.0003636942715673674

# Row 6: M2 -> ECs
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Ptprc"
$ws.Cells.Item(6, 3).Value = "Cd22"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 266.1765593333333
$ws.Cells.Item(6, 8).Value = 798.529678
$ws.Cells.Item(6, 9).Value = 0.5642668278730386
$ws.Cells.Item(6, 10).Value = 0.5642668278730386
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.5751623333333333
$ws.Cells.Item(6, 14).Value = 1.725487
$ws.Cells.Item(6, 15).Value = 0.08167574423707133
$ws.Cells.Item(6, 16).Value = 0.08167574423707133
$ws.Cells.Item(6, 17).Value = 153.0947309447984
$ws.Cells.Item(6, 18).Value = 1377.852578503186
$ws.Cells.Item(6, 19).Value = 0.04608691311482185
$ws.Cells.Item(6, 20).Value = 0.04608691311482185

# Row 7: M2 -> M2
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Ptprc"
$ws.Cells.Item(7, 3).Value = "Cd22"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 266.1765593333333
$ws.Cells.Item(7, 8).Value = 798.529678
$ws.Cells.Item(7, 9).Value = 0.5642668278730386
$ws.Cells.Item(7, 10).Value = 0.5642668278730386
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 6.466858999999999
$ws.Cells.Item(7, 14).Value = 19.400577
$ws.Cells.Item(7, 15).Value = 0.9183242557629286
$ws.Cells.Item(7, 16).Value = 0.9183242557629286
$ws.Cells.Item(7, 17).Value = 1721.3262783138
$ws.Cells.Item(7, 18).Value = 15491.93650482421
$ws.Cells.Item(7, 19).Value = 0.5181799147582168
$ws.Cells.Item(7, 20).Value = 0.5181799147582168

# Row 8: sCs -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ptprc"
$ws.Cells.Item(8, 3).Value = "Cd22"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2.0
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.044676
$ws.Cells.Item(8, 8).Value = 0.134028
$ws.Cells.Item(8, 9).Value = 0.00009470850801135487
$ws.Cells.Item(8, 10).Value = 0.00009470850801135488
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.5751623333333333
$ws.Cells.Item(8, 14).Value = 1.725487
$ws.Cells.Item(8, 15).Value = 0.08167574423707133
$ws.Cells.Item(8, 16).Value = 0.08167574423707133
$ws.Cells.Item(8, 17).Value = 0.025695952404
$ws.Cells.Item(8, 18).Value = 0.231263571636
$ws.Cells.Item(8, 19).Value = 0.00000773538787741004
$ws.Cells.Item(8, 20).Value = 0.000007735387877410042

# Row 9: sCs -> M2
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ptprc"
$ws.Cells.Item(9, 3).Value = "Cd22"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 2.0
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.044676
$ws.Cells.Item(9, 8).Value = 0.134028
$ws.Cells.Item(9, 9).Value = 0.00009470850801135487
$ws.Cells.Item(9, 10).Value = 0.00009470850801135488
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 6.466858999999999
$ws.Cells.Item(9, 14).Value = 19.400577
$ws.Cells.Item(9, 15).Value = 0.9183242557629286
$ws.Cells.Item(9, 16).Value = 0.9183242557629286
$ws.Cells.Item(9, 17).Value = 0.288913392684
$ws.Cells.Item(9, 18).Value = 2.600220534156
$ws.Cells.Item(9, 19).Value = 0.00008697312013394482
$ws.Cells.Item(9, 20).Value = 0.00008697312013394483

